$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new rows to the master-reg_center_user table, continuing the
# existing regcntr_id/usr_id sequence (rows 22-30), with the same
# lang_code/is_active/cr_by/cr_dtimes values used throughout the sheet.
$newRows = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $regcntr_id = $newRows[$i][0]
    $usr_id = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $regcntr_id
    $ws.Cells.Item($r, 2).Value = $usr_id
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Leave the selection where the author's cursor ended up: the first empty
# row below the appended data, selecting the remaining (empty) rows.
[void]$ws.Range("A31:XFD1048576").Select()
